$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as TEXT in this sheet (values like "567.69" are
# stored as text, not numbers). When a replacement price string still
# parses as a plain number, force the cell to Text format first so Excel
# doesn't auto-convert it to a numeric value on assignment.

$ws.Range('D2').Value = '68.850.80'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '3.491.74'
$ws.Range('E3').Value = '  -2.07%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.88'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.52'
$ws.Range('E6').Value = '  -2.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.614'
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('D8').Value = '3.489.24'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.643'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.84'
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000300'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.40'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '4.049.00'
$ws.Range('E15').Value = '  -2.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.18'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').Value = '68.740.33'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').Value = '3.479.82'
$ws.Range('E18').Value = '  -2.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.25'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.120'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '539.47'
$ws.Range('E21').Value = '  +13.87%  '
$ws.Range('E22').Value = '  -2.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.23'
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.98'
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.38'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.09'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.73'
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.05'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.46'
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.15'
$ws.Range('E31').Value = '  -7.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.47'
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.50'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('E34').Value = '  -4.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '568.44'
$ws.Range('E35').Value = '  -3.24%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '37.68'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.395'
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  +5.52%  '
$ws.Range('D40').Value = '0.0₃0762'
$ws.Range('E40').Value = '  -4.05%  '
$ws.Range('E41').Value = '  -4.05%  '
$ws.Range('D44').Value = '3.234.44'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.96'
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.44'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0435'
$ws.Range('E47').Value = '  -0.84%  '
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.95'
$ws.Range('E49').Value = '  -5.46%  '
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '138.54'
$ws.Range('E51').Value = '  +2.58%  '

# Rows 42/43: Kaspa and Stacks swap places (rank order changed) with updated values
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.31'
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.131'
$ws.Range('E43').Value = '  -5.11%  '

